# Apply the requested edit to "Planejamento operacional.xlsx"
# 1. Shift the dates in J/L columns (rows 19-27) back by 7 days.
# 2. Update the frozen-pane view so the top-left visible cell is A2
#    (instead of A17) and clear the stored active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- 1. Shift dates back by 7 days for the affected rows ---
$rows = @(19, 20, 21, 24, 25, 26, 27)
foreach ($r in $rows) {
    $jCell = $ws.Cells.Item($r, 10)   # column J
    $lCell = $ws.Cells.Item($r, 12)   # column L
    $jCell.Value2 = $jCell.Value2 - 7
    $lCell.Value2 = $lCell.Value2 - 7
}

# --- 2. Update the frozen pane / view and selection ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
